$values = @(
    'name',
    'DAMAGED IN AREA OF',
    'THREADS ',
    'CORRODED  at/on  ',
    'WORN BEYOND LIMITS IN THE ',
    ' LEAKING FROM',
    '? SOLENOID INOPERATIVE',
    ' ELECTRO-HYDROLIC SERVO VALVE INOPERATIVE',
    ' SHIMMED INCORRECTLY ',
    ' ELECTRO-HYDROLIC SERVO VALVE INOPERATIVE',
    ' OUT OF ADJUSTMENT',
    ' CONTAMINATED WITH',
    ' CYLINDER',
    ' SEAL (S)',
    ' BEARING (S)',
    ' PISTON CHROME',
    ' PISTON',
    ' ROD',
    ' CONNECTOR DAMAGED',
    ' UNIT OPERATING NOISY DUE TO ',
    ' WIRE (s)  DAMAGED IN AREA OF',
    ' BEARING JOURNAL WORN BEYOND LIMITS',
    ' BEARING LINER WORN BEYOND LIMITS',
    ' AXLE ASSY (BEARING JOURNALS)',
    ' AXLE ASSY (PISTON BORE)',
    ' UPPER TORQUE LINK ASSY',
    ' LOWER TORQUE LINK ASSY',
    ' MAIN FTTING',
    ' TRUNNION PINS',
    ' PISTON CROSS PIN HOLES',
    ' INSTALLATION BORE OF THE BUSHING ',
    ' UNIT COULD NOT EXTEND AND LOCK',
    ' THE PAINT FINISH - ',
    ' L.V.D.T - INOPERATIVE',
    ' LIVE TIME IS EXPIRED ',
    ' UPON ARRIVAL',
    ' LOCKWIRE DAMAGE (MISSING) FROME ',
    ' CORRODED at/on',
    ' SPRING (S)'
)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Populate Sheet2 column A with the header + list of individual defect phrases
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 1
    $ws2.Cells.Item($row, 1).Value = $values[$i]
}

# Remove the stale width formatting that used to live on column B
$ws2.Columns.Item(2).ColumnWidth = $ws2.Columns.Item(3).ColumnWidth
$ws2.Columns.Item(2).ClearFormats() | Out-Null

# Column A on Sheet2: best-fit width
$ws2.Columns.Item(1).EntireColumn.AutoFit() | Out-Null

# Update selections to match final saved state
$ws2.Activate()
$ws2.Columns.Item(1).EntireColumn.Select() | Out-Null

$ws1.Activate()
$ws1.Range("E31:F60").Select() | Out-Null
